# Synchronisation of the project's models with the google sheet metadata version 9489
# Reporting MSDS Class name upgrade in related properties
#
# The "MSDS" worksheet held the MaterialSafetyDataSheet property headers
# (materialSafetyContact, physicalChemicalProperties, ...). The "MaterialSafetyDataSheet"
# worksheet itself was empty. This change moves the header row from "MSDS" onto
# "MaterialSafetyDataSheet" and removes the now-redundant "MSDS" worksheet.

$wb = $excel.ActiveWorkbook

$msds = $wb.Worksheets.Item("MSDS")
$target = $wb.Worksheets.Item("MaterialSafetyDataSheet")

# Copy the header row (A1:O1) values from "MSDS" to "MaterialSafetyDataSheet"
$lastCol = $msds.Cells.Item(1, 1).End(-4161).Column  # xlToRight = -4161, but row has no gaps so End works; fallback below
$usedRange = $msds.UsedRange
$colCount = $usedRange.Columns.Count

for ($c = 1; $c -le $colCount; $c++) {
    $value = $msds.Cells.Item(1, $c).Value
    $target.Cells.Item(1, $c).Value = $value
}

# Remove the now-redundant "MSDS" worksheet
$msds.Delete()
